$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F8").Value = 2
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 0.97
$ws.Range("K8").Value = 0.93
$ws.Range("R8").Value = 3
$ws.Range("S8").Value = 1
$ws.Range("T8").Value = 1
$ws.Range("U8").Value = 100
$ws.Range("V8").Value = 100
